$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4) for new columns F:J
# NOTE: J4 ("...gaunt_2p") is written before I4 ("...gaunt_2p32") so the
# shared-string table is interned in the same order as the target file
# (si index 8 = "...gaunt_2p", si index 9 = "...gaunt_2p32").
$ws.Range("F4").Value = "eomip_cl_embedded_CVS_daug-acv3z_x2cmmf_gaunt_1s"
$ws.Range("G4").Value = "eomip_cl_embedded_CVS_daug-acv3z_x2cmmf_gaunt_2s"
$ws.Range("H4").Value = "eomip_cl_embedded_CVS_daug-acv3z_x2cmmf_gaunt_2p12"
$ws.Range("J4").Value = "eomip_cl_embedded_CVS_daug-acv3z_x2cmmf_gaunt_2p"
$ws.Range("I4").Value = "eomip_cl_embedded_CVS_daug-acv3z_x2cmmf_gaunt_2p32"

# Data rows 5:29 for columns F:J
$ws.Range("F5").Value = 2832.0124890000002
$ws.Range("G5").Value = 275.49195450000002
$ws.Range("H5").Value = 204.80767019999999
$ws.Range("I5").Value = 203.18072069999999
$ws.Range("J5").Value = 203.9941954
$ws.Range("F6").Value = 2831.5478170000001
$ws.Range("G6").Value = 275.06382719999999
$ws.Range("H6").Value = 204.37678320000001
$ws.Range("I6").Value = 202.7533205
$ws.Range("J6").Value = 203.56505179999999
$ws.Range("F7").Value = 2831.896009
$ws.Range("G7").Value = 275.41510249999999
$ws.Range("H7").Value = 204.7279327
$ws.Range("I7").Value = 203.09876629999999
$ws.Range("J7").Value = 203.91334950000001
$ws.Range("F8").Value = 2827.3525289999998
$ws.Range("G8").Value = 270.96737660000002
$ws.Range("H8").Value = 200.30614270000001
$ws.Range("I8").Value = 198.66106529999999
$ws.Range("J8").Value = 199.48360400000001
$ws.Range("F9").Value = 2831.924407
$ws.Range("G9").Value = 275.43348129999998
$ws.Range("H9").Value = 204.7470892
$ws.Range("I9").Value = 203.11719059999999
$ws.Range("J9").Value = 203.93213990000001
$ws.Range("F10").Value = 2830.4229989999999
$ws.Range("G10").Value = 273.84823790000002
$ws.Range("H10").Value = 203.20658689999999
$ws.Range("I10").Value = 201.57732390000001
$ws.Range("J10").Value = 202.3919554
$ws.Range("F11").Value = 2831.3227010000001
$ws.Range("G11").Value = 274.83568339999999
$ws.Range("H11").Value = 204.14910330000001
$ws.Range("I11").Value = 202.52006639999999
$ws.Range("J11").Value = 203.33458479999999
$ws.Range("F12").Value = 2831.4557909999999
$ws.Range("G12").Value = 274.93494889999999
$ws.Range("H12").Value = 204.25093910000001
$ws.Range("I12").Value = 202.62222389999999
$ws.Range("J12").Value = 203.43658149999999
$ws.Range("F13").Value = 2831.5729900000001
$ws.Range("G13").Value = 275.07081970000002
$ws.Range("H13").Value = 204.3849726
$ws.Range("I13").Value = 202.75595200000001
$ws.Range("J13").Value = 203.5704623
$ws.Range("F14").Value = 2830.7254929999999
$ws.Range("G14").Value = 274.16441090000001
$ws.Range("H14").Value = 203.52131449999999
$ws.Range("I14").Value = 201.89396500000001
$ws.Range("J14").Value = 202.70763969999999
$ws.Range("F15").Value = 2830.6114830000001
$ws.Range("G15").Value = 274.03482100000002
$ws.Range("H15").Value = 203.3928617
$ws.Range("I15").Value = 201.76548529999999
$ws.Range("J15").Value = 202.5791735
$ws.Range("F16").Value = 2830.579268
$ws.Range("G16").Value = 274.01234010000002
$ws.Range("H16").Value = 203.36952199999999
$ws.Range("I16").Value = 201.742569
$ws.Range("J16").Value = 202.55604550000001
$ws.Range("F17").Value = 2830.8898730000001
$ws.Range("G17").Value = 274.362978
$ws.Range("H17").Value = 203.71742800000001
$ws.Range("I17").Value = 202.0939218
$ws.Range("J17").Value = 202.90567490000001
$ws.Range("F18").Value = 2831.6784080000002
$ws.Range("G18").Value = 275.17331860000002
$ws.Range("H18").Value = 204.48796870000001
$ws.Range("I18").Value = 202.859814
$ws.Range("J18").Value = 203.67389130000001
$ws.Range("F19").Value = 2831.9243350000002
$ws.Range("G19").Value = 275.44120679999997
$ws.Range("H19").Value = 204.75414040000001
$ws.Range("I19").Value = 203.12455499999999
$ws.Range("J19").Value = 203.93934770000001
$ws.Range("F20").Value = 2831.1869689999999
$ws.Range("G20").Value = 274.6309329
$ws.Range("H20").Value = 203.9877774
$ws.Range("I20").Value = 202.36316020000001
$ws.Range("J20").Value = 203.1754688
$ws.Range("F21").Value = 2830.836135
$ws.Range("G21").Value = 274.27420510000002
$ws.Range("H21").Value = 203.6309718
$ws.Range("I21").Value = 202.00571400000001
$ws.Range("J21").Value = 202.8183429
$ws.Range("F22").Value = 2829.7700100000002
$ws.Range("G22").Value = 273.16596800000002
$ws.Range("H22").Value = 202.52674999999999
$ws.Range("I22").Value = 200.90221679999999
$ws.Range("J22").Value = 201.71448340000001
$ws.Range("F23").Value = 2831.916522
$ws.Range("G23").Value = 275.42743369999999
$ws.Range("H23").Value = 204.74069220000001
$ws.Range("I23").Value = 203.11369619999999
$ws.Range("J23").Value = 203.9271942
$ws.Range("F24").Value = 2832.3043010000001
$ws.Range("G24").Value = 275.83231599999999
$ws.Range("H24").Value = 205.14413579999999
$ws.Range("I24").Value = 203.51640399999999
$ws.Range("J24").Value = 204.33026989999999
$ws.Range("F25").Value = 2830.5625530000002
$ws.Range("G25").Value = 273.98376389999999
$ws.Range("H25").Value = 203.3418666
$ws.Range("I25").Value = 201.712918
$ws.Range("J25").Value = 202.5273923
$ws.Range("F26").Value = 2829.776758
$ws.Range("G26").Value = 273.19024810000002
$ws.Range("H26").Value = 202.54930630000001
$ws.Range("I26").Value = 200.9238618
$ws.Range("J26").Value = 201.73658399999999
$ws.Range("F27").Value = 2831.6812319999999
$ws.Range("G27").Value = 275.20270729999999
$ws.Range("H27").Value = 204.51525860000001
$ws.Range("I27").Value = 202.88700019999999
$ws.Range("J27").Value = 203.70112940000001
$ws.Range("F28").Value = 2831.314069
$ws.Range("G28").Value = 274.85109360000001
$ws.Range("H28").Value = 204.16245240000001
$ws.Range("I28").Value = 202.540402
$ws.Range("J28").Value = 203.35142719999999
$ws.Range("F29").Value = 2831.6636859999999
$ws.Range("G29").Value = 275.16677650000003
$ws.Range("H29").Value = 204.4808165
$ws.Range("I29").Value = 202.85215220000001
$ws.Range("J29").Value = 203.6664844

# Update selection to match post-edit state (K4)
$ws.Range("K4").Select()
